$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21) for every data row
# (rows 2 through 482).
$ws.Range("C2:C482").Value = 45190

# Row 4 (record "A 31224-2022") gained one extra signal species
# ("Tallfingersvamp"), so the related counters and the species list
# text must be updated as well.
$ws.Range("I4").Value = 1
$ws.Range("Q4").Value = 5

$nl = "`r`n"
$ws.Range("R4").Value = "Knärot" + $nl + "Orange taggsvamp" + $nl + "Svart taggsvamp" + $nl + "Talltaggsvamp" + $nl + "Tallfingersvamp"
